# validatie.ipynb werkt nu vanaf de omgezette karteringen
#
# This script updates the validationinfo worksheet:
#  - renames header A1 "naam_in_overzicht_vegkart" -> "path_pred"
#  - renames header D1 "path" -> "path_true"
#  - replaces column A (rows 2-9) with the real (converted) source shapefile
#    paths instead of a duplicate of the "naam_in_overzicht_habkart" name
#  - adds two new columns: G "werkt_nu" (boolean) and H "reden dat hij niet werkt"
#  - marks every kartering as working (TRUE) except the RuitenAa2020 one,
#    which gets FALSE plus an explanation in column H
#  - removes the grey highlight fill that was on A6
#  - widens the new columns G/H and moves the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("A1").Value = "path_pred"
$ws.Range("D1").Value = "path_true"
$ws.Range("G1").Value = "werkt_nu"
$ws.Range("H1").Value = "reden dat hij niet werkt"

# --- New "path_pred" values for the Groningen karteringen (rows 2-9) ---
$ws.Range("A2").Value = "./GR/SBB Westerwolde 2020/GISbestanden/vlakken.shp"
$ws.Range("A3").Value = "./GR/SBB ZWK 2010/0814_Tussen de Gasten 2010/ZWK0814_2010.shp"
$ws.Range("A4").Value = "./GR/SBB ZWK 2010/0815_Trimunt_2010/vlakken.shp"
$ws.Range("A5").Value = "./GR/NM vegetatiekartering RuitenAa2020 edited/vegkart_RuitenA_2020/vegkart_RuitenA_2020.shp"
$ws.Range("A6").Value = "./GR/SGL Zuidlaardermeer 2019/kartering zuidlaardermeer 2019 definitief/GIS-bestanden Zuidlaardermeer 2019/GIS_Vlakken_zuidlaardermeer_2019/vlakken.shp"
$ws.Range("A7").Value = "./GR/SGL Hunzedal en Leekstermeer2021/2021 Vegetatiekartering Leekstermeer2021/GIS bestanden Onlanden 2021/Vegetatiekartering_Leekstermeer2021.shp"
$ws.Range("A8").Value = "./GR/SBB Lauwersmeer 2015/vlakken.shp"
$ws.Range("A9").Value = "./GR/SGL Hunzedal en Leekstermeer2021/2021 Vegetatie- en plantensoortenkartering Hunzedal concept/gis/Vegetatiekartering_Hunzedal2021.shp"

# remove the grey fill that used to flag row 6 as a mismatch
$ws.Range("A6").ClearFormats()

# --- "werkt_nu" / "reden dat hij niet werkt" columns --------------------
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 7).Value = $true
}
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = "De kartering die wij hebben heebt geen unieke ElmID kolom; dus kunnen vlakken niet met zekerheid aan vegetatietypen gekoppeld worden"

# --- Column widths for the two new columns ------------------------------
$ws.Columns.Item(7).ColumnWidth = 10.666666666666666
$ws.Columns.Item(8).ColumnWidth = 164.5

# --- View housekeeping ---------------------------------------------------
[void]$ws.Range("C41").Select()

Write-Host "done"
